$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1007.2727
$ws.Range("J40").Value = 1157
$ws.Range("L40").Value = 1157
$ws.Range("N40").Value = -1507

$ws.Range("H64").Value = 2981.4546
$ws.Range("J64").Value = 2999.25
$ws.Range("L64").Value = 2999.25
$ws.Range("N64").Value = -3495.25

$ws.Range("H67").Value = 2981.4546
$ws.Range("J67").Value = 2999.25
$ws.Range("L67").Value = 2999.25
$ws.Range("N67").Value = -4715.25

$ws.Range("H115").Value = 5554.857
$ws.Range("I115").Value = 5576.8
$ws.Range("J115").Value = 5500
$ws.Range("K115").Value = 16730.4
$ws.Range("L115").Value = 16500
$ws.Range("M115").Value = -15163.4
$ws.Range("N115").Value = -19634

$ws.Range("H127").Value = 768.1177
$ws.Range("J127").Value = 1194.4286
$ws.Range("L127").Value = 3583.2858
$ws.Range("N127").Value = -13503.2858

$ws.Range("H137").Value = 6947780.5
$ws.Range("I137").Value = 13892098
$ws.Range("J137").Value = 3463.3333
$ws.Range("K137").Value = 41676294
$ws.Range("L137").Value = 10389.9999
$ws.Range("M137").Value = -41673744
$ws.Range("N137").Value = -15489.9999

$ws.Range("H138").Value = 4816.9375
$ws.Range("I138").Value = 6813.4287
$ws.Range("J138").Value = 4571.7544
$ws.Range("K138").Value = 20440.2861
$ws.Range("L138").Value = 13715.2632
$ws.Range("M138").Value = -15300.2861
$ws.Range("N138").Value = -23995.2632

$ws.Range("H140").Value = 76544
$ws.Range("J140").Value = 76544
$ws.Range("L140").Value = 76544
$ws.Range("N140").Value = -86904

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 29999
$ws.Range("J23").Value = 29999
$ws.Range("L23").Value = 29999
$ws.Range("N23").Value = -30517

$ws.Range("H45").Value = 2151.7
$ws.Range("I45").Value = 1481.7778
$ws.Range("J45").Value = 2699.818
$ws.Range("K45").Value = 1481.7778
$ws.Range("L45").Value = 2699.818
$ws.Range("M45").Value = -1104.7778
$ws.Range("N45").Value = -3453.818

$ws.Range("H97").Value = 53368.58
$ws.Range("I97").Value = 56294
$ws.Range("J97").Value = 711
$ws.Range("K97").Value = 56294
$ws.Range("L97").Value = 711
$ws.Range("M97").Value = -55798
$ws.Range("N97").Value = -1703

$ws.Range("H110").Value = 750
$ws.Range("I110").Value = 595
$ws.Range("J110").Value = 905
$ws.Range("K110").Value = 595
$ws.Range("L110").Value = 905
$ws.Range("M110").Value = 1450
$ws.Range("N110").Value = -4995

$ws.Range("H122").Value = 101082
$ws.Range("I122").Value = 121047
$ws.Range("J122").Value = 1257
$ws.Range("K122").Value = 363141
$ws.Range("L122").Value = 3771
$ws.Range("M122").Value = -360691
$ws.Range("N122").Value = -8671

$ws.Range("H132").Value = 2657574.2
$ws.Range("I132").Value = 5918.933
$ws.Range("J132").Value = 5498633.5
$ws.Range("K132").Value = 17756.799
$ws.Range("L132").Value = 16495900.5
$ws.Range("M132").Value = -15226.799
$ws.Range("N132").Value = -16500960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 167215
$ws.Range("I94").Value = 500245
$ws.Range("K94").Value = 500245
$ws.Range("M94").Value = -499794

$ws.Range("H134").Value = 2341.6562
$ws.Range("I134").Value = 2262.7036
$ws.Range("J134").Value = 2768
$ws.Range("K134").Value = 6788.110799999999
$ws.Range("L134").Value = 8304
$ws.Range("M134").Value = -4253.110799999999
$ws.Range("N134").Value = -13374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 11463.789
$ws.Range("J4").Value = 11463.789
$ws.Range("L4").Value = 11463.789
$ws.Range("N4").Value = -11687.789

$ws.Range("H31").Value = 9083.754000000001
$ws.Range("I31").Value = 1101.2142
$ws.Range("J31").Value = 11682.721
$ws.Range("K31").Value = 1101.2142
$ws.Range("L31").Value = 11682.721
$ws.Range("M31").Value = -806.2141999999999
$ws.Range("N31").Value = -12272.721

$ws.Range("H34").Value = 9083.754000000001
$ws.Range("I34").Value = 1101.2142
$ws.Range("J34").Value = 11682.721
$ws.Range("K34").Value = 1101.2142
$ws.Range("L34").Value = 11682.721
$ws.Range("M34").Value = -899.2141999999999
$ws.Range("N34").Value = -12086.721

$ws.Range("H122").Value = 2227.611
$ws.Range("I122").Value = 1616.6
$ws.Range("J122").Value = 2462.6155
$ws.Range("K122").Value = 4849.799999999999
$ws.Range("L122").Value = 7387.8465
$ws.Range("M122").Value = -2399.799999999999
$ws.Range("N122").Value = -12287.8465

$ws.Range("H132").Value = 18019412
$ws.Range("I132").Value = 19231952
$ws.Range("J132").Value = 15153406
$ws.Range("K132").Value = 57695856
$ws.Range("L132").Value = 45460218
$ws.Range("M132").Value = -57693326
$ws.Range("N132").Value = -45465278

$ws.Range("H134").Value = 3523826
$ws.Range("I134").Value = 3908960.5
$ws.Range("J134").Value = 2597.5715
$ws.Range("K134").Value = 11726881.5
$ws.Range("L134").Value = 7792.7145
$ws.Range("M134").Value = -11724346.5
$ws.Range("N134").Value = -12862.7145

$ws.Range("H141").Value = 92954.19500000001
$ws.Range("J141").Value = 91768.414
$ws.Range("L141").Value = 91768.414
$ws.Range("N141").Value = -102128.414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1027.7179
$ws.Range("I5").Value = 782.04614
$ws.Range("J5").Value = 2256.077
$ws.Range("K5").Value = 2346.13842
$ws.Range("L5").Value = 6768.231000000001
$ws.Range("M5").Value = -2234.13842
$ws.Range("N5").Value = -6992.231000000001

$ws.Range("H12").Value = 59.954544
$ws.Range("I12").Value = 40.933334
$ws.Range("J12").Value = 100.71429
$ws.Range("K12").Value = 122.800002
$ws.Range("L12").Value = 302.14287
$ws.Range("M12").Value = 50.19999799999999
$ws.Range("N12").Value = -648.14287

$ws.Range("H44").Value = 406
$ws.Range("I44").Value = 367.6
$ws.Range("J44").Value = 502
$ws.Range("K44").Value = 1102.8
$ws.Range("L44").Value = 1506
$ws.Range("M44").Value = -704.8000000000002
$ws.Range("N44").Value = -2302

$ws.Range("H107").Value = 18519472
$ws.Range("J107").Value = 27779042
$ws.Range("L107").Value = 83337126
$ws.Range("N107").Value = -83340966

$ws.Range("H113").Value = 913.2941
$ws.Range("I113").Value = 761.3570999999999
$ws.Range("J113").Value = 1158.7307
$ws.Range("K113").Value = 2284.0713
$ws.Range("L113").Value = 3476.1921
$ws.Range("M113").Value = -114.0712999999996
$ws.Range("N113").Value = -7816.1921

$ws.Range("H121").Value = 731.1111
$ws.Range("I121").Value = 316.66666
$ws.Range("J121").Value = 1560
$ws.Range("K121").Value = 949.9999799999999
$ws.Range("L121").Value = 4680
$ws.Range("M121").Value = 360.0000200000001
$ws.Range("N121").Value = -7300

$ws.Range("H131").Value = 3891.814
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 3891.814
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 11675.442
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -21755.442

$ws.Range("H133").Value = 16715.723
$ws.Range("I133").Value = 7990
$ws.Range("J133").Value = 23696.3
$ws.Range("K133").Value = 23970
$ws.Range("L133").Value = 71088.89999999999
$ws.Range("M133").Value = -18910
$ws.Range("N133").Value = -81208.89999999999

$ws.Range("H135").Value = 1027.7179
$ws.Range("I135").Value = 782.04614
$ws.Range("J135").Value = 2256.077
$ws.Range("K135").Value = 7038.415260000001
$ws.Range("L135").Value = 20304.693
$ws.Range("M135").Value = -4503.415260000001
$ws.Range("N135").Value = -25374.693

$ws.Range("H136").Value = 3653.6875
$ws.Range("I136").Value = 1650.6923
$ws.Range("J136").Value = 12333.333
$ws.Range("K136").Value = 4952.0769
$ws.Range("L136").Value = 36999.999
$ws.Range("M136").Value = 147.9231
$ws.Range("N136").Value = -47199.999

$ws.Range("H139").Value = 336343.16
$ws.Range("I139").Value = 386034.7
$ws.Range("J139").Value = 13348.25
$ws.Range("K139").Value = 1158104.1
$ws.Range("L139").Value = 40044.75
$ws.Range("M139").Value = -1152964.1
$ws.Range("N139").Value = -50324.75

$ws.Range("H141").Value = 7333.857
$ws.Range("I141").Value = 4380
$ws.Range("J141").Value = 10287.714
$ws.Range("K141").Value = 13140
$ws.Range("L141").Value = 30863.142
$ws.Range("M141").Value = -7960
$ws.Range("N141").Value = -41223.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 146602.86
$ws.Range("I97").Value = 146602.86
$ws.Range("K97").Value = 146602.86
$ws.Range("M97").Value = -146106.86

$ws.Range("H132").Value = 38468708
$ws.Range("I132").Value = 83344370
$ws.Range("J132").Value = 3853.7144
$ws.Range("K132").Value = 250033110
$ws.Range("L132").Value = 11561.1432
$ws.Range("M132").Value = -250030580
$ws.Range("N132").Value = -16621.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9992
$ws.Range("J2").Value = 9992
$ws.Range("L2").Value = 9992
$ws.Range("N2").Value = -10216

$ws.Range("H132").Value = 3332.1082
$ws.Range("I132").Value = 2854.6667
$ws.Range("J132").Value = 3784.4211
$ws.Range("K132").Value = 8564.000100000001
$ws.Range("L132").Value = 11353.2633
$ws.Range("M132").Value = -6034.000100000001
$ws.Range("N132").Value = -16413.2633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 70015
$ws.Range("J22").Value = 70015
$ws.Range("L22").Value = 70015
$ws.Range("N22").Value = -70601

$ws.Range("H51").Value = 17535
$ws.Range("I51").Value = 20070
$ws.Range("K51").Value = 20070
$ws.Range("M51").Value = -19560

$ws.Range("H126").Value = 1969.5555
$ws.Range("I126").Value = 1871
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 5613
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -3143
$ws.Range("N126").Value = -11440.0001

$ws.Range("H132").Value = 8336044
$ws.Range("I132").Value = 2683.16
$ws.Range("K132").Value = 8049.48
$ws.Range("M132").Value = -5519.48

$ws.Range("H136").Value = 5519.304
$ws.Range("I136").Value = 6199.353
$ws.Range("J136").Value = 3592.5
$ws.Range("K136").Value = 18598.059
$ws.Range("L136").Value = 10777.5
$ws.Range("M136").Value = -16048.059
$ws.Range("N136").Value = -15877.5
